# Commit: "added the functionality of checking for US numbers with area code +1"
#
# The sample data in column A is replaced/extended with a new set of Nigerian
# (and a few +234) phone-number strings used to exercise the new "US number"
# validation logic, growing the sheet from A1:A10 to A1:A28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-parsed as a number by Excel
# (losing a leading zero or the leading "+"), so force them to text format first.
$numericTextCells = @("A1", "A5", "A9", "A13", "A17", "A21", "A28")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Populate A1:A28 with the new phone-number sample data.
$ws.Range("A1").Value = "07025664511"
$ws.Range("A2").Value = "0702 5664 511"
$ws.Range("A3").Value = "0703 331 3456"
$ws.Range("A4").Value = "070 3444 2345"
$ws.Range("A5").Value = "07125664511"
$ws.Range("A6").Value = "0712 5664 511"
$ws.Range("A7").Value = "0713 331 3456"
$ws.Range("A8").Value = "071 3444 2345"
$ws.Range("A9").Value = "08025664511"
$ws.Range("A10").Value = "0802 5664 511"
$ws.Range("A11").Value = "0803 331 3456"
$ws.Range("A12").Value = "080 3444 2345"
$ws.Range("A13").Value = "08125664511"
$ws.Range("A14").Value = "0812 5664 511"
$ws.Range("A15").Value = "0813 331 3456"
$ws.Range("A16").Value = "081 3444 2345"
$ws.Range("A17").Value = "09125664511"
$ws.Range("A18").Value = "0912 5664 511"
$ws.Range("A19").Value = "0913 331 3456"
$ws.Range("A20").Value = "091 3444 2345"
$ws.Range("A21").Value = "09025664511"
$ws.Range("A22").Value = "0902 5664 511"
$ws.Range("A23").Value = "0903 331 3456"
$ws.Range("A24").Value = "090 3444 2345"
$ws.Range("A25").Value = "+234 445 334 5555"
$ws.Range("A26").Value = "+234 445 3344 555"
$ws.Range("A27").Value = "+234 4453 345 555"
$ws.Range("A28").Value = "+2348123456780"

Write-Output ("Dimension now spans rows 1 to 28; A1=" + $ws.Range("A1").Value2 + " A28=" + $ws.Range("A28").Value2)
